$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15 (shifts GBDT..XGBoost rows down by one)
$ws.Rows.Item(15).Insert()

# Populate the new row 15 with the DeepCNN model name; leave hyperparameter/value
# columns blank (they still carry the inlineStr type in the target workbook).
$ws.Range("A15").Value = "DeepCNN"
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = ""
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("F15").Value = ""
